$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 254, shifting existing rows 254-265 down to 255-266
$ws.Rows.Item(254).Insert()

# Populate the newly inserted row 254 with data (copy pattern from the row below / surrounding rows)
$ws.Cells.Item(254, 1).Value = 10
$ws.Cells.Item(254, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(254, 3).Value = "La Araucanía"
$ws.Cells.Item(254, 4).Value = 45265
$ws.Cells.Item(254, 5).Value = 9
$ws.Cells.Item(254, 6).Value = "Fruta"
$ws.Cells.Item(254, 7).Value = 100107
$ws.Cells.Item(254, 8).Value = "Otros"
$ws.Cells.Item(254, 9).Value = 100107002
$ws.Cells.Item(254, 10).Value = "Chirimoya"
$ws.Cells.Item(254, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(254, 12).Value = "Primera"
$ws.Cells.Item(254, 13).Value = 90
$ws.Cells.Item(254, 14).Value = 2300
$ws.Cells.Item(254, 15).Value = 2300
$ws.Cells.Item(254, 16).Value = 2300
$ws.Cells.Item(254, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(254, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(254, 19).Value = 2300
$ws.Cells.Item(254, 20).Value = 1

# Ensure the date cell (column D) keeps the same number format / style as the rest of the column
$ws.Cells.Item(254, 4).NumberFormat = $ws.Cells.Item(255, 4).NumberFormat
